# Apply updated crypto price/volume snapshot values from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (prevents Excel from auto-converting
# numeric-looking strings like "1.00" or "61.863.00" into numbers/dates).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
}

# Row 2
Set-TextValue 'D2' '61.863.00'
Set-TextValue 'E2' '  -2.87%  '
# Row 3
Set-TextValue 'D3' '2.493.55'
Set-TextValue 'E3' '  -5.16%  '
# Row 4
Set-TextValue 'E4' '  +0.05%  '
# Row 5
Set-TextValue 'D5' '555.42'
Set-TextValue 'E5' '  -3.66%  '
# Row 6
Set-TextValue 'D6' '147.56'
Set-TextValue 'E6' '  -4.75%  '
# Row 7
Set-TextValue 'E7' '  -0.02%  '
# Row 8
Set-TextValue 'E8' '  -3.14%  '
# Row 9
Set-TextValue 'D9' '2.489.93'
# Row 10
Set-TextValue 'D10' '0.109'
Set-TextValue 'E10' '  -7.49%  '
# Row 11
Set-TextValue 'E11' '  -6.24%  '
# Row 12
Set-TextValue 'E12' '  -1.38%  '
# Row 13
Set-TextValue 'D13' '0.362'
Set-TextValue 'E13' '  -5.18%  '
# Row 14
Set-TextValue 'D14' '26.45'
Set-TextValue 'E14' '  -6.81%  '
# Row 15
Set-TextValue 'D15' '2.939.18'
Set-TextValue 'E15' '  -5.25%  '
# Row 16
Set-TextValue 'E16' '  -7.16%  '
# Row 17
Set-TextValue 'D17' '61.737.76'
Set-TextValue 'E17' '  -2.97%  '
# Row 18
Set-TextValue 'D18' '2.501.54'
Set-TextValue 'E18' '  -5.02%  '
# Row 19
Set-TextValue 'D19' '11.22'
Set-TextValue 'E19' '  -7.40%  '
# Row 20
Set-TextValue 'D20' '7.04'
Set-TextValue 'E20' '  -7.27%  '
# Row 21
Set-TextValue 'D21' '4.23'
Set-TextValue 'E21' '  -6.62%  '
# Row 22
Set-TextValue 'D22' '324.17'
Set-TextValue 'E22' '  -5.99%  '
# Row 23
Set-TextValue 'E23' '  +0.00%  '
# Row 24
Set-TextValue 'E24' '  -4.84%  '
# Row 25
Set-TextValue 'D25' '64.29'
Set-TextValue 'E25' '  -5.41%  '
# Row 26
Set-TextValue 'D26' '0.0000101'
Set-TextValue 'E26' '  -6.65%  '
# Row 27
Set-TextValue 'B27' 'Fetch.AI'
Set-TextValue 'C27' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D27' '1.55'
Set-TextValue 'E27' '  -3.53%  '
# Row 28
Set-TextValue 'B28' 'WrappedeETH'
Set-TextValue 'C28' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D28' '2.610.81'
Set-TextValue 'E28' '  -5.07%  '
# Row 29
Set-TextValue 'B29' 'Bittensor'
Set-TextValue 'C29' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D29' '540.26'
Set-TextValue 'E29' '  -10.99%  '
# Row 30
Set-TextValue 'B30' 'Binance-PegBSC-USD'
Set-TextValue 'C30' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.01%  '
# Row 31
Set-TextValue 'B31' 'InternetComputer(DFINITY)'
Set-TextValue 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D31' '8.44'
Set-TextValue 'E31' '  -8.72%  '
# Row 32
Set-TextValue 'D32' '7.60'
Set-TextValue 'E32' '  -4.61%  '
# Row 33
Set-TextValue 'E33' '  -5.21%  '
# Row 34
Set-TextValue 'E34' '  -6.80%  '
# Row 35
Set-TextValue 'D35' '1.61'
Set-TextValue 'E35' '  -7.56%  '
# Row 36
Set-TextValue 'D36' '6.01'
Set-TextValue 'E36' '  -9.23%  '
# Row 37
Set-TextValue 'D37' '4.95'
Set-TextValue 'E37' '  -8.20%  '
# Row 39
Set-TextValue 'E39' '  -4.04%  '
# Row 40
Set-TextValue 'D40' '18.64'
Set-TextValue 'E40' '  -5.51%  '
# Row 41
Set-TextValue 'D41' '148.49'
Set-TextValue 'E41' '  -1.06%  '
# Row 43
Set-TextValue 'D43' '1.00'
Set-TextValue 'E43' '  +0.09%  '
# Row 44
Set-TextValue 'D44' '40.42'
Set-TextValue 'E44' '  -3.08%  '
# Row 45
Set-TextValue 'D45' '2.37'
Set-TextValue 'E45' '  -6.30%  '
# Row 46
Set-TextValue 'D46' '149.38'
Set-TextValue 'E46' '  -5.99%  '
# Row 47
Set-TextValue 'D47' '3.65'
Set-TextValue 'E47' '  -6.35%  '
# Row 48
Set-TextValue 'D48' '21.23'
Set-TextValue 'E48' '  -14.52%  '
# Row 49
Set-TextValue 'D49' '0.0539'
Set-TextValue 'E49' '  -7.99%  '
# Row 50
Set-TextValue 'D50' '0.601'
Set-TextValue 'E50' '  -4.67%  '
# Row 51
Set-TextValue 'D51' '0.0951'
Set-TextValue 'E51' '  -4.61%  '
